# Weekly update: a new "Apio" price record (week of 2021-11-xx, serial 44518)
# was inserted into the daily logic subset sheet ahead of the existing
# (now shifted down) rows, per the "Fruta / hortaliza, semanal" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 98, pushing the existing rows 98..147 down to 99..148.
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(98, 1).Value = 7
$ws.Cells.Item(98, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(98, 3).Value = "Ñuble"
$ws.Cells.Item(98, 4).Value = 44518
$ws.Cells.Item(98, 5).Value = 16
$ws.Cells.Item(98, 6).Value = 100112017
$ws.Cells.Item(98, 7).Value = "Apio"
$ws.Cells.Item(98, 8).Value = "Americana (o)"
$ws.Cells.Item(98, 9).Value = "Primera"
$ws.Cells.Item(98, 10).Value = 100
$ws.Cells.Item(98, 11).Value = 8000
$ws.Cells.Item(98, 12).Value = 9000
$ws.Cells.Item(98, 13).Value = 8500
$ws.Cells.Item(98, 14).Value = "`$/docena de matas"
$ws.Cells.Item(98, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(98, 16).Value = 1417
$ws.Cells.Item(98, 17).Value = 6
$ws.Cells.Item(98, 18).Value = "Hortaliza"
